$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the erroneous "licensed drivers" figure for 1936 (row 7): 755 -> 755765
$ws.Range("B7").Value = 755765

# Append historical crash-statistics rows for years 1964-2010 (sheet rows 35-81)
$newRows = @(
    @(1964,2694023,111232,1424,54560),
    @(1965,2739138,128462,1611,60917),
    @(1966,2821648,139781,1596,65210),
    @(1967,3004654,145008,1719,67280),
    @(1968,3128509,155127,1586,71520),
    @(1969,3247979,169395,1683,74902),
    @(1970,3422892,141609,1535,75126),
    @(1971,3563197,158831,1769,84650),
    @(1972,3688541,189494,1934,95181),
    @(1973,3841628,193021,1959,97790),
    @(1974,3972980,204271,1748,98673),
    @(1975,4160623,213689,1800,97034),
    @(1976,4315925,211865,1511,83736),
    @(1977,4562903,218567,1420,95964),
    @(1978,4725546,186363,1450,94979),
    @(1979,4858351,197196,1560,101321),
    @(1980,4993531,196501,1508,101367),
    @(1981,5123177,198372,1445,100321),
    @(1982,5247198,187943,1138,92815),
    @(1983,5380259,181999,1204,91706),
    @(1984,5513911,194782,1132,97230),
    @(1985,5660422,189750,1191,109169),
    @(1986,5817799,187286,1102,108839),
    @(1987,5978105,203431,1229,121089),
    @(1988,6118112,228398,1237,118158),
    @(1989,6290424,247038,1286,120652),
    @(1990,6448883,220188,1120,101575),
    @(1991,6574231,213669,1102,90519),
    @(1992,6688761,224249,1090,91025),
    @(1993,6823428,228834,1135,91149),
    @(1994,6983960,226996,999,90030),
    @(1995,7086018,219085,999,89572),
    @(1996,7258167,215024,929,88445),
    @(1997,7537607,221500,899,85527),
    @(1998,7727756,213356,854,83154),
    @(1999,7918314,221962,868,84062),
    @(2000,8121374,240630,849,85009),
    @(2001,8266616,234004,845,81782),
    @(2002,8413504,244642,873,84192),
    @(2003,8541555,246463,831,77879),
    @(2004,8655597,231548,799,73008),
    @(2005,8762210,230258,766,71850),
    @(2006,8867965,216247,769,68793),
    @(2007,8945397,233487,765,67175),
    @(2008,9042286,229196,631,62743),
    @(2009,9101938,216315,564,62562),
    @(2010,9245267,215533,579,64514)
)

$startRow = 35
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Match the post-edit selection / scrolled view from the source workbook
$ws.Range("C74").Select()
$excel.ActiveWindow.TopLeftCell = $ws.Range("A64")
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.ScrollColumn = 1
